$wb = $excel.ActiveWorkbook

# --- Sheet "öffentliche Portale": update the Bund URL (B5) and drop its hyperlink ---
$ws = $wb.Worksheets.Item("öffentliche Portale")
$ws.Activate()

# Update cell value to shortened URL
$ws.Range("B5").Value = "https://www.service.bund.de/"

# Remove the hyperlink attached to B5 (keeps the cell text, drops the link)
if ($ws.Hyperlinks.Count -gt 0) {
    foreach ($hl in @($ws.Hyperlinks)) {
        if ($hl.Range.Row -eq 5 -and $hl.Range.Column -eq 2) {
            $hl.Delete()
        }
    }
}

# Move the active selection on this sheet to B8
$ws.Range("B8").Select()

# --- Workbook view: mark window as minimized ---
$wb.Windows.Item(1).WindowState = -4140  # xlMinimized
